$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.467036
$ws.Cells.Item(2, 8).Value = 1.401108
$ws.Cells.Item(2, 9).Value = 0.001972893265924874
$ws.Cells.Item(2, 10).Value = 0.001972893265924874
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.19897
$ws.Cells.Item(2, 14).Value = 0.5969100000000001
$ws.Cells.Item(2, 15).Value = 0.001481973067923264
$ws.Cells.Item(2, 16).Value = 0.001481973067923264
$ws.Cells.Item(2, 17).Value = 0.09292615292
$ws.Cells.Item(2, 18).Value = 0.8363353762800001
$ws.Cells.Item(2, 19).Value = [double]"2.923774685987834E-06"
$ws.Cells.Item(2, 20).Value = [double]"2.923774685987834E-06"
$ws.Cells.Item(3, 7).Value = 0.467036
$ws.Cells.Item(3, 8).Value = 1.401108
$ws.Cells.Item(3, 9).Value = 0.001972893265924874
$ws.Cells.Item(3, 10).Value = 0.001972893265924874
$ws.Cells.Item(3, 15).Value = 0.001601566712998508
$ws.Cells.Item(3, 16).Value = 0.001601566712998507
$ws.Cells.Item(3, 17).Value = 0.1004251942933333
$ws.Cells.Item(3, 18).Value = 0.90382674864
$ws.Cells.Item(3, 19).Value = [double]"3.159720183004192E-06"
$ws.Cells.Item(3, 20).Value = [double]"3.159720183004192E-06"
$ws.Cells.Item(4, 7).Value = 0.467036
$ws.Cells.Item(4, 8).Value = 1.401108
$ws.Cells.Item(4, 9).Value = 0.001972893265924874
$ws.Cells.Item(4, 10).Value = 0.001972893265924874
$ws.Cells.Item(4, 13).Value = 0.245373
$ws.Cells.Item(4, 14).Value = 0.736119
$ws.Cells.Item(4, 15).Value = 0.001827592991885888
$ws.Cells.Item(4, 16).Value = 0.001827592991885888
$ws.Cells.Item(4, 17).Value = 0.114598024428
$ws.Cells.Item(4, 18).Value = 1.031382219852
$ws.Cells.Item(4, 19).Value = [double]"3.605645906543161E-06"
$ws.Cells.Item(4, 20).Value = [double]"3.605645906543161E-06"
$ws.Cells.Item(5, 7).Value = 0.467036
$ws.Cells.Item(5, 8).Value = 1.401108
$ws.Cells.Item(5, 9).Value = 0.001972893265924874
$ws.Cells.Item(5, 10).Value = 0.001972893265924874
$ws.Cells.Item(5, 13).Value = 133.60083
$ws.Cells.Item(5, 14).Value = 400.80249
$ws.Cells.Item(5, 15).Value = 0.9950888672271924
$ws.Cells.Item(5, 16).Value = 0.9950888672271923
$ws.Cells.Item(5, 17).Value = 62.39639723988
$ws.Cells.Item(5, 18).Value = 561.56757515892
$ws.Cells.Item(5, 19).Value = 0.001963204125149339
$ws.Cells.Item(5, 20).Value = 0.001963204125149339
$ws.Cells.Item(6, 9).Value = 0.001374344438283074
$ws.Cells.Item(6, 10).Value = 0.001374344438283074
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.19897
$ws.Cells.Item(6, 14).Value = 0.5969100000000001
$ws.Cells.Item(6, 15).Value = 0.001481973067923264
$ws.Cells.Item(6, 16).Value = 0.001481973067923264
$ws.Cells.Item(6, 17).Value = 0.06473362935666667
$ws.Cells.Item(6, 18).Value = 0.5826026642100001
$ws.Cells.Item(6, 19).Value = [double]"2.036741443585642E-06"
$ws.Cells.Item(6, 20).Value = [double]"2.036741443585642E-06"
$ws.Cells.Item(7, 9).Value = 0.001374344438283074
$ws.Cells.Item(7, 10).Value = 0.001374344438283074
$ws.Cells.Item(7, 15).Value = 0.001601566712998508
$ws.Cells.Item(7, 16).Value = 0.001601566712998507
$ws.Cells.Item(7, 19).Value = [double]"2.201104304548803E-06"
$ws.Cells.Item(7, 20).Value = [double]"2.201104304548803E-06"
$ws.Cells.Item(8, 9).Value = 0.001374344438283074
$ws.Cells.Item(8, 10).Value = 0.001374344438283074
$ws.Cells.Item(8, 13).Value = 0.245373
$ws.Cells.Item(8, 14).Value = 0.736119
$ws.Cells.Item(8, 15).Value = 0.001827592991885888
$ws.Cells.Item(8, 16).Value = 0.001827592991885888
$ws.Cells.Item(8, 17).Value = 0.07983055152099999
$ws.Cells.Item(8, 18).Value = 0.718474963689
$ws.Cells.Item(8, 19).Value = [double]"2.511742263843492E-06"
$ws.Cells.Item(8, 20).Value = [double]"2.511742263843492E-06"
$ws.Cells.Item(9, 9).Value = 0.001374344438283074
$ws.Cells.Item(9, 10).Value = 0.001374344438283074
$ws.Cells.Item(9, 13).Value = 133.60083
$ws.Cells.Item(9, 14).Value = 400.80249
$ws.Cells.Item(9, 15).Value = 0.9950888672271924
$ws.Cells.Item(9, 16).Value = 0.9950888672271923
$ws.Cells.Item(9, 17).Value = 43.46618390191
$ws.Cells.Item(9, 18).Value = 391.19565511719
$ws.Cells.Item(9, 19).Value = 0.001367594850271096
$ws.Cells.Item(9, 20).Value = 0.001367594850271096
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.02089333333333333
$ws.Cells.Item(10, 8).Value = 0.06268
$ws.Cells.Item(10, 9).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(10, 10).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.19897
$ws.Cells.Item(10, 14).Value = 0.5969100000000001
$ws.Cells.Item(10, 15).Value = 0.001481973067923264
$ws.Cells.Item(10, 16).Value = 0.001481973067923264
$ws.Cells.Item(10, 17).Value = 0.004157146533333334
$ws.Cells.Item(10, 18).Value = 0.03741431880000001
$ws.Cells.Item(10, 19).Value = [double]"1.307980521970594E-07"
$ws.Cells.Item(10, 20).Value = [double]"1.307980521970594E-07"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.02089333333333333
$ws.Cells.Item(11, 8).Value = 0.06268
$ws.Cells.Item(11, 9).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(11, 10).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(11, 15).Value = 0.001601566712998508
$ws.Cells.Item(11, 16).Value = 0.001601566712998507
$ws.Cells.Item(11, 17).Value = 0.004492623822222222
$ws.Cells.Item(11, 18).Value = 0.0404336144
$ws.Cells.Item(11, 19).Value = [double]"1.413533154265786E-07"
$ws.Cells.Item(11, 20).Value = [double]"1.413533154265786E-07"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.02089333333333333
$ws.Cells.Item(12, 8).Value = 0.06268
$ws.Cells.Item(12, 9).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(12, 10).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(12, 13).Value = 0.245373
$ws.Cells.Item(12, 14).Value = 0.736119
$ws.Cells.Item(12, 15).Value = 0.001827592991885888
$ws.Cells.Item(12, 16).Value = 0.001827592991885888
$ws.Cells.Item(12, 17).Value = 0.00512665988
$ws.Cells.Item(12, 18).Value = 0.04613993892
$ws.Cells.Item(12, 19).Value = [double]"1.613022589422981E-07"
$ws.Cells.Item(12, 20).Value = [double]"1.613022589422981E-07"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.02089333333333333
$ws.Cells.Item(13, 8).Value = 0.06268
$ws.Cells.Item(13, 9).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(13, 10).Value = [double]"8.825939892440207E-05"
$ws.Cells.Item(13, 13).Value = 133.60083
$ws.Cells.Item(13, 14).Value = 400.80249
$ws.Cells.Item(13, 15).Value = 0.9950888672271924
$ws.Cells.Item(13, 16).Value = 0.9950888672271923
$ws.Cells.Item(13, 17).Value = 2.7913666748
$ws.Cells.Item(13, 18).Value = 25.1223000732
$ws.Cells.Item(13, 19).Value = [double]"8.782594529783614E-05"
$ws.Cells.Item(13, 20).Value = [double]"8.782594529783612E-05"
$ws.Cells.Item(14, 7).Value = 235.9131673333333
$ws.Cells.Item(14, 8).Value = 707.739502
$ws.Cells.Item(14, 9).Value = 0.9965645028968676
$ws.Cells.Item(14, 10).Value = 0.9965645028968676
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.19897
$ws.Cells.Item(14, 14).Value = 0.5969100000000001
$ws.Cells.Item(14, 15).Value = 0.001481973067923264
$ws.Cells.Item(14, 16).Value = 0.001481973067923264
$ws.Cells.Item(14, 17).Value = 46.93964290431334
$ws.Cells.Item(14, 18).Value = 422.4567861388201
$ws.Cells.Item(14, 19).Value = 0.001476881753741493
$ws.Cells.Item(14, 20).Value = 0.001476881753741493
$ws.Cells.Item(15, 7).Value = 235.9131673333333
$ws.Cells.Item(15, 8).Value = 707.739502
$ws.Cells.Item(15, 9).Value = 0.9965645028968676
$ws.Cells.Item(15, 10).Value = 0.9965645028968676
$ws.Cells.Item(15, 15).Value = 0.001601566712998508
$ws.Cells.Item(15, 16).Value = 0.001601566712998507
$ws.Cells.Item(15, 17).Value = 50.72762199446223
$ws.Cells.Item(15, 18).Value = 456.54859795016
$ws.Cells.Item(15, 19).Value = 0.001596064535195528
$ws.Cells.Item(15, 20).Value = 0.001596064535195528
$ws.Cells.Item(16, 7).Value = 235.9131673333333
$ws.Cells.Item(16, 8).Value = 707.739502
$ws.Cells.Item(16, 9).Value = 0.9965645028968676
$ws.Cells.Item(16, 10).Value = 0.9965645028968676
$ws.Cells.Item(16, 13).Value = 0.245373
$ws.Cells.Item(16, 14).Value = 0.736119
$ws.Cells.Item(16, 15).Value = 0.001827592991885888
$ws.Cells.Item(16, 16).Value = 0.001827592991885888
$ws.Cells.Item(16, 17).Value = 57.886721608082
$ws.Cells.Item(16, 18).Value = 520.980494472738
$ws.Cells.Item(16, 19).Value = 0.001821314301456558
$ws.Cells.Item(16, 20).Value = 0.001821314301456558
$ws.Cells.Item(17, 7).Value = 235.9131673333333
$ws.Cells.Item(17, 8).Value = 707.739502
$ws.Cells.Item(17, 9).Value = 0.9965645028968676
$ws.Cells.Item(17, 10).Value = 0.9965645028968676
$ws.Cells.Item(17, 13).Value = 133.60083
$ws.Cells.Item(17, 14).Value = 400.80249
$ws.Cells.Item(17, 15).Value = 0.9950888672271924
$ws.Cells.Item(17, 16).Value = 0.9950888672271923
$ws.Cells.Item(17, 17).Value = 31518.19496366222
$ws.Cells.Item(17, 18).Value = 283663.75467296
$ws.Cells.Item(17, 19).Value = 0.991670242306474
$ws.Cells.Item(17, 20).Value = 0.991670242306474
